$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - AK10 unchanged ("test"); AL10:AO10 -> "testing"
$ws.Range("AL10:AO10").Value = "testing"

# Row 11 - AK11:AO11 -> "testing"
$ws.Range("AK11:AO11").Value = "testing"

# Row 14 - AK14:AO14 -> "test"
$ws.Range("AK14:AO14").Value = "test"

# Row 15 - AK15 -> "testing"; AL15:AO15 -> "test"
$ws.Range("AK15").Value = "testing"
$ws.Range("AL15:AO15").Value = "test"

# Row 16 - AK16 -> "test"; AL16 -> "t"
$ws.Range("AK16").Value = "test"
$ws.Range("AL16").Value = "t"

# Row 17 - AK17 -> "test"; AL17 -> "te"
$ws.Range("AK17").Value = "test"
$ws.Range("AL17").Value = "te"
